$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 38, shifting all rows 38:65 down
# to 39:66 (preserving their data/format), then populate the newly
# inserted row 38 with the new weekly record.
$ws.Rows.Item(38).Insert()

$ws.Range("A38").Value = 1
$ws.Range("B38").Value = 'Agrícola del Norte S.A. de Arica'
$ws.Range("C38").Value = 'Arica y Parinacota'
$ws.Range("D38").Value2 = 44586
$ws.Range("E38").Value = 15
$ws.Range("F38").Value = 100112038
$ws.Range("G38").Value = 'Cebollín baby'
$ws.Range("H38").Value = 'Sin especificar'
$ws.Range("I38").Value = 'Primera'
$ws.Range("J38").Value = 200
$ws.Range("K38").Value = 5000
$ws.Range("L38").Value = 5500
$ws.Range("M38").Value = 5250
$ws.Range("N38").Value = '$/paquete 1,5 a 2 kilos'
$ws.Range("O38").Value = 'Región de Arica y Parinacota'
$ws.Range("P38").Value = 2625
$ws.Range("Q38").Value = 2
$ws.Range("R38").Value = 'Hortaliza'
